$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 202
$ws.Range("F4").Value = 679
$ws.Range("F6").Value = 2202
$ws.Range("F7").Value = 1304
$ws.Range("F8").Value = 774
$ws.Range("F10").Value = 18
$ws.Range("F11").Value = 2790
$ws.Range("F12").Value = 19
$ws.Range("F14").Value = 1067
$ws.Range("F15").Value = 560
$ws.Range("F17").Value = 871
$ws.Range("F18").Value = 72
$ws.Range("F19").Value = 77
$ws.Range("F21").Value = 93
$ws.Range("F22").Value = 597
$ws.Range("F23").Value = 579
$ws.Range("F24").Value = 257
$ws.Range("F26").Value = 951
$ws.Range("F27").Value = 4848
$ws.Range("F28").Value = 353
$ws.Range("F29").Value = 124
$ws.Range("F30").Value = 53

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 353
$ws.Range("F21").Value = 23
$ws.Range("F22").Value = 290
$ws.Range("F23").Value = 32
$ws.Range("F24").Value = 29
$ws.Range("F25").Value = 336
$ws.Range("F27").Value = 517
$ws.Range("F31").Value = 46
$ws.Range("F37").Value = 680

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 1494
$ws.Range("F5").Value = 604
$ws.Range("F6").Value = 359
$ws.Range("F7").Value = 331

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1494
$ws.Range("F4").Value = 604
$ws.Range("F5").Value = 202
$ws.Range("F6").Value = 359
$ws.Range("F9").Value = 679
$ws.Range("F10").Value = 353
$ws.Range("F13").Value = 2202
$ws.Range("F14").Value = 1304
$ws.Range("F15").Value = 774
$ws.Range("F19").Value = 18
$ws.Range("F20").Value = 2790
$ws.Range("F21").Value = 19
$ws.Range("F24").Value = 1067
$ws.Range("F25").Value = 560
$ws.Range("F27").Value = 331
$ws.Range("F29").Value = 871
$ws.Range("F30").Value = 871
$ws.Range("F31").Value = 72
$ws.Range("F32").Value = 23
$ws.Range("F33").Value = 290
$ws.Range("F34").Value = 77
$ws.Range("F35").Value = 93
$ws.Range("F36").Value = 32
$ws.Range("F37").Value = 29
$ws.Range("F38").Value = 597
$ws.Range("F39").Value = 579
$ws.Range("F40").Value = 336
$ws.Range("F41").Value = 517
$ws.Range("F42").Value = 257
$ws.Range("F45").Value = 951
$ws.Range("F46").Value = 4848
$ws.Range("F47").Value = 46
$ws.Range("F48").Value = 353
$ws.Range("F49").Value = 124
$ws.Range("F50").Value = 680

